$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 19:47:50"
$ws1.Range("A3").Value = "Total filas: 339"
$ws1.Range("C141").Value = "23_HERNANDEZ"
$ws1.Range("C142").Value = "27_EL RETIRO"
$ws1.Range("A184").Value = "14:32:44"
$ws1.Range("C184").Value = "14X44_ABASTO"
$ws1.Range("D184").Value = 1
$ws1.Range("A185").Value = "13:55:43"
$ws1.Range("C185").Value = "215C_EL PATO"
$ws1.Range("D185").Value = 38
$ws1.Range("A204").Value = "13:55:43"
$ws1.Range("C204").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D204").Value = 118
$ws1.Range("C205").Value = "15X38_ABASTO"
$ws1.Range("A206").Value = "14:32:44"
$ws1.Range("C206").Value = "10_OLMOS"
$ws1.Range("D206").Value = 81
$ws1.Range("C269").Value = "15_ABASTO"
$ws1.Range("C270").Value = "16_P MOR-SANTA ANA"
$ws1.Range("A301").Value = "17:35:41"
$ws1.Range("C301").Value = "27_EL RETIRO"
$ws1.Range("D301").Value = 101
$ws1.Range("A302").Value = "17:55:25"
$ws1.Range("C302").Value = "17_ROMERO"
$ws1.Range("D302").Value = 81
$ws1.Range("A329").Value = "19:47:50"
$ws1.Range("B329").Value = "20:33"
$ws1.Range("D329").Value = 46
$ws1.Range("A330").Value = "19:35:34"
$ws1.Range("B330").Value = "20:34"
$ws1.Range("C330").Value = "16_SANTA ANA"
$ws1.Range("D330").Value = 59
$ws1.Range("A331").Value = "19:11:44"
$ws1.Range("B331").Value = "20:41"
$ws1.Range("D331").Value = 90
$ws1.Range("A332").Value = "18:52:29"
$ws1.Range("B332").Value = "20:42"
$ws1.Range("D332").Value = 110
$ws1.Range("A333").Value = "19:35:34"
$ws1.Range("B333").Value = "20:43"
$ws1.Range("C333").Value = "17_ROMERO"
$ws1.Range("D333").Value = 68
$ws1.Range("A334").Value = "19:47:50"
$ws1.Range("B334").Value = "20:45"
$ws1.Range("C334").Value = "17_ROMERO"
$ws1.Range("D334").Value = 58
$ws1.Range("A335").Value = "18:52:29"
$ws1.Range("B335").Value = "20:47"
$ws1.Range("C335").Value = "215B_EL PATO"
$ws1.Range("D335").Value = 115
$ws1.Range("A336").Value = "19:35:34"
$ws1.Range("B336").Value = "20:55"
$ws1.Range("C336").Value = "23_HERNANDEZ"
$ws1.Range("D336").Value = 80
$ws1.Range("A337").Value = "19:11:44"
$ws1.Range("B337").Value = "20:56"
$ws1.Range("C337").Value = "27_EL RETIRO"
$ws1.Range("D337").Value = 105
$ws1.Range("A338").Value = "19:11:44"
$ws1.Range("B338").Value = "21:06"
$ws1.Range("C338").Value = "10_OLMOS"
$ws1.Range("D338").Value = 115
$ws1.Range("A339").Value = "19:47:50"
$ws1.Range("B339").Value = "21:09"
$ws1.Range("C339").Value = "15_ABASTO"
$ws1.Range("D339").Value = 82
$ws1.Range("A340").Value = "19:35:34"
$ws1.Range("B340").Value = "21:10"
$ws1.Range("C340").Value = "15_ABASTO"
$ws1.Range("D340").Value = 95
$ws1.Range("E340").Value = "LP1912"
$ws1.Range("A341").Value = "19:35:34"
$ws1.Range("B341").Value = "21:28"
$ws1.Range("C341").Value = "11_ETCHEVERRY"
$ws1.Range("D341").Value = 113
$ws1.Range("E341").Value = "LP1912"
$ws1.Range("A342").Value = "19:47:50"
$ws1.Range("B342").Value = "21:33"
$ws1.Range("C342").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D342").Value = 106
$ws1.Range("E342").Value = "LP1912"
$ws1.Range("A343").Value = "19:35:34"
$ws1.Range("B343").Value = "21:34"
$ws1.Range("C343").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D343").Value = 119
$ws1.Range("E343").Value = "LP1912"
$ws1.Range("A344").Value = "19:47:50"
$ws1.Range("B344").Value = "21:45"
$ws1.Range("C344").Value = "14X44_ABASTO"
$ws1.Range("D344").Value = 118
$ws1.Range("E344").Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 19:47:50"

# --- Sheet: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 19:47:50"
$ws3.Range("A3").Value = "Total filas: 48"
$ws3.Range("A52").Value = "19:47:50"
$ws3.Range("B52").Value = "21:27"
$ws3.Range("D52").Value = 100
$ws3.Range("A53").Value = "19:35:34"
$ws3.Range("B53").Value = "21:30"
$ws3.Range("C53").Value = "215C_LA PLATA"
$ws3.Range("D53").Value = 115
$ws3.Range("E53").Value = "L6203"
